{"js": "// 1. Remove the existing \"_GoBack\" bookmark (it sat right after\n//    \"...we used python and \" in the SQL/Python paragraph).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst body = context.document.body;\n\n// 2. Rewrite the \"significant words\" sentence:\n//      \"...ficant words from each decade were discovered  from among the\n//      words that were  not used in previous decades. This analysis\n//      allowed us to see emerging themes in Terror over time.\"\n//    becomes\n//      \"...ficant words were discovered  and were filtered to not include\n//      significant words from previous decades. This analysis allowed us\n//      to see emerging themes in terror over time.\"\nconst oldSentence1 =\n  \"ficant words from each decade were discovered  from among the words that were  not used in previous decades. This analysis allowed us to see emerging themes in Terror over time.\";\nconst newSentence1 =\n  \"ficant words were discovered  and were filtered to not include significant words from previous decades. This analysis allowed us to see emerging themes in terror over time.\";\n\nconst results1 = body.search(oldSentence1, { matchCase: true });\nresults1.load(\"items\");\nawait context.sync();\nif (results1.items.length === 0) {\n  throw new Error(\"Could not locate the 'significant words' sentence to replace.\");\n}\nresults1.items[0].insertText(newSentence1, \"Replace\");\nawait context.sync();\n\n// 3. Rewrite the Saladin / Muslim-military-campaign sentence: add an \"s\"\n//    to turn \"campaign\" into \"campaigns\", and rephrase \"He is referenced\n//    in the current era associated...\" into \"His reference in the\n//    current era is associated...\".\nconst oldSentence2 =\n  \" Muslim military campaign against the Crusader states in the Levant. He is referenced in the current era associated to the brutality of the methods by which he gained and retained power.\";\nconst newSentence2 =\n  \" Muslim military campaigns against the Crusader states in the Levant. His reference in the current era is associated to the brutality of the methods by which he gained and retained power.\";\n\nconst results2 = body.search(oldSentence2, { matchCase: true });\nresults2.load(\"items\");\nawait context.sync();\nif (results2.items.length === 0) {\n  throw new Error(\"Could not locate the 'Muslim military campaign' sentence to replace.\");\n}\nresults2.items[0].insertText(newSentence2, \"Replace\");\nawait context.sync();\n\n// 4. Re-insert the \"_GoBack\" bookmark right after the new \"campaigns\" word\n//    (i.e. just before \" against the Crusader states...\").\nconst results3 = body.search(\"Muslim military campaigns\", { matchCase: true });\nresults3.load(\"items\");\nawait context.sync();\nif (results3.items.length === 0) {\n  throw new Error(\"Could not locate 'Muslim military campaigns' to anchor the bookmark.\");\n}\nconst campaignsRange = results3.items[0];\nconst afterCampaigns = campaignsRange.getRange(\"End\");\nafterCampaigns.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Remove the existing \"_GoBack\" bookmark (it sat right after\n#    \"...we used python and \" in the SQL/Python paragraph).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2. Rewrite the \"significant words\" sentence.\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"ficant words from each decade were discovered  from among the words that were  not used in previous decades. This analysis allowed us to see emerging themes in Terror over time.\"\n$find1.Replacement.Text = \"ficant words were discovered  and were filtered to not include significant words from previous decades. This analysis allowed us to see emerging themes in terror over time.\"\n$find1.MatchCase = $true\n$find1.Execute($null, $null, $null, $null, $null, $null, $true, $null, $null, $null, 2)\n\n# 3. Rewrite the Saladin/Muslim military campaign sentence (adds an \"s\" to\n#    \"campaign\" -> \"campaigns\" and rephrases the closing sentence).\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \" Muslim military campaign against the Crusader states in the Levant. He is referenced in the current era associated to the brutality of the methods by which he gained and retained power.\"\n$find2.Replacement.Text = \" Muslim military campaigns against the Crusader states in the Levant. His reference in the current era is associated to the brutality of the methods by which he gained and retained power.\"\n$find2.MatchCase = $true\n$find2.Execute($null, $null, $null, $null, $null, $null, $true, $null, $null, $null, 2)\n\n# 4. Re-insert the \"_GoBack\" bookmark right after the new \"campaigns\" word.\n$range = $d.Content\n$find3 = $range.Find\n$find3.ClearFormatting()\n$find3.Text = \"Muslim military campaigns\"\n$find3.MatchCase = $true\n$find3.Execute() | Out-Null\n$range.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $range)\n"}
